# Apply the Ghostbusters Spectral Search content edits.
$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title / heading text loses the " ᐈ Demo Game" suffix (occurs twice: H1 + bold run near bottom)
Replace-AllText "Play Ghostbusters Spectral Search Free ᐈ Demo Game" "Play Ghostbusters Spectral Search Free"

# "What we like" bullet list
Replace-AllText "High RTP percentage of 96.4%" "High RTP percentage"
Replace-AllText "Unique gameplay experience using a high-tech detection device" "Complex gameplay style"
Replace-AllText "Various locations to search for ghosts" "Variety of locations to search for ghosts"
Replace-AllText "Varying symbol values for payout" "Unique and thrilling gameplay experience"

# "What we don't like" bullet list
Replace-AllText "Gameplay may be too complex for some players" "May take some time to understand the gameplay"
Replace-AllText "Limited instant payouts" "Limited symbol values compared to other slot games"

# Meta description / italic summary paragraph
Replace-AllText "Get ready to play Ghostbusters Spectral Search for free! Read our review to discover the pros and cons of this high RTP online slot game and try it now." "Read a review of Ghostbusters Spectral Search, a free slot game with unique gameplay."
